$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-12-31 Wednesday"; new = "2026-01-01 Thursday"},
    @{old = "624÷4=156, 0"; new = "768÷5=153, 3"},
    @{old = "268÷9=29, 7"; new = "690÷9=76, 6"},
    @{old = "681÷4=170, 1"; new = "350÷9=38, 8"},
    @{old = "112÷4=28, 0"; new = "240÷7=34, 2"},
    @{old = "702÷6=117, 0"; new = "355÷7=50, 5"},
    @{old = "516÷9=57, 3"; new = "527÷2=263, 1"},
    @{old = "229÷8=28, 5"; new = "268÷5=53, 3"},
    @{old = "665÷7=95, 0"; new = "274÷7=39, 1"},
    @{old = "467÷7=66, 5"; new = "777÷3=259, 0"},
    @{old = "537÷4=134, 1"; new = "116÷2=58, 0"},
    @{old = "226÷5=45, 1"; new = "874÷8=109, 2"},
    @{old = "875÷7=125, 0"; new = "188÷5=37, 3"},
    @{old = "483÷2=241, 1"; new = "438÷4=109, 2"},
    @{old = "578÷2=289, 0"; new = "793÷9=88, 1"},
    @{old = "947÷9=105, 2"; new = "660÷6=110, 0"},
    @{old = "223÷2=111, 1"; new = "850÷8=106, 2"},
    @{old = "794÷7=113, 3"; new = "821÷4=205, 1"},
    @{old = "163÷9=18, 1"; new = "968÷4=242, 0"},
    @{old = "725÷6=120, 5"; new = "880÷5=176, 0"},
    @{old = "319÷7=45, 4"; new = "659÷8=82, 3"},
    @{old = "805÷9=89, 4"; new = "589÷7=84, 1"},
    @{old = "421÷7=60, 1"; new = "195÷7=27, 6"},
    @{old = "283÷4=70, 3"; new = "956÷6=159, 2"},
    @{old = "439÷3=146, 1"; new = "843÷8=105, 3"},
    @{old = "798÷4=199, 2"; new = "848÷2=424, 0"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
